$wb = $excel.ActiveWorkbook

# Sheet "展览" — update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 234
$ws1.Range("F4").Value = 856
$ws1.Range("F5").Value = 76

# Sheet "全部类型" — update "想去人数" (F column) values (rows shifted by +1)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 234
$ws4.Range("F5").Value = 856
$ws4.Range("F6").Value = 76
